$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6
$ws.Range("F4").Value = 5
$ws.Range("F6").Value = -5
$ws.Range("F12").Value = -5
